$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that needs to move from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2-271).
$ws.Range("C2:C271").Value = 45188
